# Add 2022-Q3 data:
#  - summary sheet ("总计") gets a new top row for 2022-Q3
#  - a new worksheet "2022-Q3" is inserted right after "总计" (before the
#    former "2022-Q2" tab), populated with the Q3 fund-holdings data. It is
#    created as a copy of the "2022-Q2" sheet (same layout/header/styles)
#    and then its three data rows are overwritten with the Q3 numbers.

$wb = $excel.ActiveWorkbook

function Set-TextCell($rng, $val) {
    # Force the cell to be stored as text (matches the source data, which
    # keeps numeric-looking figures like "9.55" / "0.5090" / "016620" as
    # strings so leading zeros / trailing zeros survive), then drop back to
    # the default "Normal" style so no stray number-format style sticks.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. Summary sheet ("总计"): insert a new row 2 for the 2022-Q3 totals.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# The inserted row inherits row-1's (header) formatting; reset B2:D2 back
# to the plain/default style used by every other data row, and give A2 the
# same "index column" style (s=2) the rest of column A already uses.
$summary.Range("B2:D2").Style = "Normal"
$summary.Range("A3").Copy() | Out-Null
$summary.Range("A2").PasteSpecial(-4122) | Out-Null

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 3
$summary.Range("D2").Value = 0.57

# Re-sequence the index column (A) for every row that shifted down so it
# stays a contiguous 0..7 counter.
for ($r = 3; $r -le 9; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------
# 2. New "2022-Q3" worksheet: copy of "2022-Q2", inserted right before it,
#    then overwrite the three fund rows with the Q3 figures.
# ---------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("2022-Q2")
$srcSheet.Copy($srcSheet)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# Row 2
Set-TextCell $newSheet.Range("B2") "519198"
$newSheet.Range("C2").Value = "万家颐和灵活配置混合A"
Set-TextCell $newSheet.Range("D2") "9.55"
Set-TextCell $newSheet.Range("E2") "93.93"
Set-TextCell $newSheet.Range("F2") "5.33"
Set-TextCell $newSheet.Range("G2") "0.5090"
$newSheet.Range("H2").Value = 9

# Row 3
Set-TextCell $newSheet.Range("B3") "519197"
$newSheet.Range("C3").Value = "万家颐达灵活配置混合"
Set-TextCell $newSheet.Range("D3") "2.25"
Set-TextCell $newSheet.Range("E3") "45.36"
Set-TextCell $newSheet.Range("F3") "2.19"
Set-TextCell $newSheet.Range("G3") "0.0493"
$newSheet.Range("H3").Value = 8

# Row 4
Set-TextCell $newSheet.Range("B4") "016620"
$newSheet.Range("C4").Value = "万家颐和灵活配置混合C"
Set-TextCell $newSheet.Range("D4") "0.18"
Set-TextCell $newSheet.Range("E4") "93.93"
Set-TextCell $newSheet.Range("F4") "5.33"
Set-TextCell $newSheet.Range("G4") "0.0096"
$newSheet.Range("H4").Value = 9
